# Insert a new daily price record as row 39 (pushing the existing rows
# 39-125 down to 40-126, i.e. one row lower than before).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = 7
$ws.Range("B39").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C39").Value = "Ñuble"
$ws.Range("D39").Value = 45044
$ws.Range("E39").Value = 16
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100108
$ws.Range("H39").Value = "Tropicales y subtropicales"
$ws.Range("I39").Value = 100108002
$ws.Range("J39").Value = "Mango"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 70
$ws.Range("N39").Value = 7000
$ws.Range("O39").Value = 8000
$ws.Range("P39").Value = 7571
$ws.Range("Q39").Value = "$/bandeja 4 kilos"
$ws.Range("R39").Value = "Perú"
$ws.Range("S39").Value = 1893
$ws.Range("T39").Value = 4
